$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1803.6
$ws.Range("I28").Value = 222.57143
$ws.Range("J28").Value = 3187
$ws.Range("K28").Value = 222.57143
$ws.Range("L28").Value = 3187
$ws.Range("M28").Value = 262.42857
$ws.Range("N28").Value = -4157

# ALC row 61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 1473
$ws.Range("I61").Value = 1473
$ws.Range("K61").Value = 4419
$ws.Range("M61").Value = -4247

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4358
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4358
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 889
$ws.Range("I132").Value = 879.28125
$ws.Range("K132").Value = 2637.84375
$ws.Range("M132").Value = -107.84375

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1964.579
$ws.Range("I137").Value = 1450.1111
$ws.Range("J137").Value = 2427.6
$ws.Range("K137").Value = 4350.3333
$ws.Range("L137").Value = 7282.799999999999
$ws.Range("M137").Value = -1800.3333
$ws.Range("N137").Value = -12382.8

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3763.9355
$ws.Range("I138").Value = 5357.0835
$ws.Range("K138").Value = 16071.2505
$ws.Range("M138").Value = -10931.2505

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2627.9248
$ws.Range("I32").Value = 2113.4458
$ws.Range("K32").Value = 2113.4458
$ws.Range("M32").Value = -1826.4458

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1347.5
$ws.Range("J97").Value = 1883.3334
$ws.Range("L97").Value = 1883.3334
$ws.Range("N97").Value = -2875.3334

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 944
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1198.75
$ws.Range("I122").Value = 895
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 2685
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -235
$ws.Range("N122").Value = -8800

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1717.9395
$ws.Range("I132").Value = 1489.0714
$ws.Range("J132").Value = 2999.6
$ws.Range("K132").Value = 4467.2142
$ws.Range("L132").Value = 8998.799999999999
$ws.Range("M132").Value = -1937.2142
$ws.Range("N132").Value = -14058.8

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13385.643
$ws.Range("I134").Value = 13385.643
$ws.Range("K134").Value = 40156.929
$ws.Range("M134").Value = -37621.929

# CRP row 57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2175252.8
$ws.Range("I58").Value = 3345742
$ws.Range("J58").Value = 1486.7142
$ws.Range("K58").Value = 3345742
$ws.Range("L58").Value = 1486.7142
$ws.Range("M58").Value = -3345539
$ws.Range("N58").Value = -1892.7142

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2576.0312
$ws.Range("I132").Value = 1887.3478
$ws.Range("K132").Value = 5662.0434
$ws.Range("M132").Value = -3132.0434

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2175252.8
$ws.Range("I136").Value = 3345742
$ws.Range("J136").Value = 1486.7142
$ws.Range("K136").Value = 10037226
$ws.Range("L136").Value = 4460.142599999999
$ws.Range("M136").Value = -10034676
$ws.Range("N136").Value = -9560.142599999999

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2314.077
$ws.Range("J68").Value = 2779.889
$ws.Range("L68").Value = 8339.667000000001
$ws.Range("N68").Value = -9961.667000000001

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2314.077
$ws.Range("J71").Value = 2779.889
$ws.Range("L71").Value = 25019.001
$ws.Range("N71").Value = -33131.001

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1570.5333
$ws.Range("I107").Value = 1184.2858
$ws.Range("J107").Value = 1641.6842
$ws.Range("K107").Value = 3552.8574
$ws.Range("L107").Value = 4925.0526
$ws.Range("M107").Value = -1632.8574
$ws.Range("N107").Value = -8765.052599999999

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2229.8333
$ws.Range("I80").Value = 2503
$ws.Range("J80").Value = 2175.2
$ws.Range("K80").Value = 2503
$ws.Range("L80").Value = 2175.2
$ws.Range("M80").Value = -1505
$ws.Range("N80").Value = -4171.2

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2229.8333
$ws.Range("I83").Value = 2503
$ws.Range("J83").Value = 2175.2
$ws.Range("K83").Value = 12515
$ws.Range("L83").Value = 10876
$ws.Range("M83").Value = -7523
$ws.Range("N83").Value = -20860

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2566839
$ws.Range("I132").Value = 2960506.5
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 8881519.5
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -8878989.5
$ws.Range("N132").Value = -29060

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7244.385
$ws.Range("I16").Value = 9110.444
$ws.Range("J16").Value = 3045.75
$ws.Range("K16").Value = 9110.444
$ws.Range("L16").Value = 3045.75
$ws.Range("M16").Value = -8940.444
$ws.Range("N16").Value = -3385.75

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1213.3636
$ws.Range("I46").Value = 475
$ws.Range("J46").Value = 1377.4445
$ws.Range("K46").Value = 475
$ws.Range("L46").Value = 1377.4445
$ws.Range("M46").Value = -287
$ws.Range("N46").Value = -1753.4445

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 209.22728
$ws.Range("I55").Value = 231.57143
$ws.Range("J55").Value = 170.125
$ws.Range("K55").Value = 231.57143
$ws.Range("L55").Value = 170.125
$ws.Range("M55").Value = -58.57142999999999
$ws.Range("N55").Value = -516.125

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3982.0435
$ws.Range("I132").Value = 1082.1666
$ws.Range("J132").Value = 5005.5293
$ws.Range("K132").Value = 3246.4998
$ws.Range("L132").Value = 15016.5879
$ws.Range("M132").Value = -716.4998000000001
$ws.Range("N132").Value = -20076.5879

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 20499.75
$ws.Range("J96").Value = 20499.75
$ws.Range("L96").Value = 20499.75
$ws.Range("N96").Value = -23245.75

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 711.6667
$ws.Range("I113").Value = 366.66666
$ws.Range("K113").Value = 1099.99998
$ws.Range("M113").Value = 1070.00002

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 118782.75

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10895065
$ws.Range("I136").Value = 26457190
$ws.Range("K136").Value = 79371570
$ws.Range("M136").Value = -79369020
